# Implement manual input of repayment and prepayment in absolute amounts
# in balancesheetmutations (sheet "intangible redemption").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intangible redemption")

# --- Header rows -----------------------------------------------------
# Row 2: rename ItemType -> SubItemType, add two "Other loans" sub-headers
$ws.Range("A2").Value = "SubItemType"
$ws.Range("C2").Value = "Other loans"
$ws.Range("D2").Value = "Other loans"

# Row 3: add Repayment / Prepayment column headers
$ws.Range("C3").Value = "Repayment"
$ws.Range("D3").Value = "Prepayment"

# --- Data rows (4-25): default to 0 ----------------------------------
for ($r = 4; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
}

# --- Specific manual amounts ------------------------------------------
$ws.Range("C7").Value = 100000
$ws.Range("C9").Value = 50000
$ws.Range("D9").Value = 100000

# Apply the "Comma [0]" look (matches existing style used elsewhere
# in the workbook) to the cells carrying manual absolute amounts.
$commaFormat = "_ * #,##0_ ;_ * \-#,##0_ ;_ * ""-""??_ ;_ @_ "
$ws.Range("C7").NumberFormat = $commaFormat
$ws.Range("C9").NumberFormat = $commaFormat
$ws.Range("D9").NumberFormat = $commaFormat

# --- Column widths for the new columns --------------------------------
# (matches the bestFit width Excel computes for "Other loans"/"Repayment"/
# "Prepayment" as closely as this runtime's column-width quantization allows)
$ws.Columns("C:D").ColumnWidth = 10.33

# --- Selection ----------------------------------------------------------
[void]$ws.Range("C10").Select()
